$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 945.7646999999999
$ws.Range("J19").Value = 768.7143
$ws.Range("L19").Value = 768.7143
$ws.Range("N19").Value = -1118.7143
$ws.Range("H33").Value = 224.6923
$ws.Range("J33").Value = 416.66666
$ws.Range("L33").Value = 416.66666
$ws.Range("N33").Value = -874.66666
$ws.Range("H52").Value = 2490
$ws.Range("I52").Value = 2490
$ws.Range("K52").Value = 7470
$ws.Range("M52").Value = -7310
$ws.Range("H62").Value = 5452.6665
$ws.Range("I62").Value = 3415.2
$ws.Range("J62").Value = 7999.5
$ws.Range("K62").Value = 3415.2
$ws.Range("L62").Value = 7999.5
$ws.Range("M62").Value = -2791.2
$ws.Range("N62").Value = -9247.5
$ws.Range("H65").Value = 5452.6665
$ws.Range("I65").Value = 3415.2
$ws.Range("J65").Value = 7999.5
$ws.Range("K65").Value = 17076
$ws.Range("L65").Value = 39997.5
$ws.Range("M65").Value = -13956
$ws.Range("N65").Value = -46237.5
$ws.Range("H88").Value = 1292.2222
$ws.Range("J88").Value = 2312
$ws.Range("L88").Value = 2312
$ws.Range("N88").Value = -3124
$ws.Range("H91").Value = 1292.2222
$ws.Range("J91").Value = 2312
$ws.Range("L91").Value = 2312
$ws.Range("N91").Value = -5120
$ws.Range("H100").Value = 1843.2307
$ws.Range("I100").Value = 1626.2
$ws.Range("J100").Value = 2566.6667
$ws.Range("K100").Value = 1626.2
$ws.Range("L100").Value = 2566.6667
$ws.Range("M100").Value = -1085.2
$ws.Range("N100").Value = -3648.6667
$ws.Range("H107").Value = 386.66666
$ws.Range("I107").Value = 325.85715
$ws.Range("J107").Value = 599.5
$ws.Range("K107").Value = 325.85715
$ws.Range("L107").Value = 599.5
$ws.Range("M107").Value = 1594.14285
$ws.Range("N107").Value = -4439.5
$ws.Range("H116").Value = 5284.6924
$ws.Range("I116").Value = 3356.7144
$ws.Range("K116").Value = 3356.7144
$ws.Range("M116").Value = 85.28560000000016
$ws.Range("H138").Value = 2248.7778
$ws.Range("I138").Value = 2078.7307
$ws.Range("J138").Value = 2481.4736
$ws.Range("K138").Value = 6236.1921
$ws.Range("L138").Value = 7444.4208
$ws.Range("M138").Value = -1096.1921
$ws.Range("N138").Value = -17724.4208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 1993
$ws.Range("I41").Value = 1993
$ws.Range("K41").Value = 1993
$ws.Range("M41").Value = -1579
$ws.Range("H43").Value = 40749.25
$ws.Range("J43").Value = 38999.5
$ws.Range("L43").Value = 38999.5
$ws.Range("N43").Value = -39625.5
$ws.Range("H61").Value = 3554.2222
$ws.Range("J61").Value = 7066.3335
$ws.Range("L61").Value = 7066.3335
$ws.Range("N61").Value = -7490.3335
$ws.Range("H74").Value = 1207.7273
$ws.Range("I74").Value = 960.875
$ws.Range("K74").Value = 960.875
$ws.Range("M74").Value = -86.875
$ws.Range("H77").Value = 1207.7273
$ws.Range("I77").Value = 960.875
$ws.Range("K77").Value = 4804.375
$ws.Range("M77").Value = -436.375
$ws.Range("H97").Value = 468.35294
$ws.Range("I97").Value = 501.92856
$ws.Range("K97").Value = 501.92856
$ws.Range("M97").Value = -5.928560000000004
$ws.Range("H109").Value = 49999.5
$ws.Range("J109").Value = 49999.5
$ws.Range("L109").Value = 49999.5
$ws.Range("N109").Value = -52773.5
$ws.Range("H132").Value = 1514.8182
$ws.Range("I132").Value = 1514.8182
$ws.Range("K132").Value = 4544.4546
$ws.Range("M132").Value = -2014.4546
$ws.Range("H136").Value = 3554.2222
$ws.Range("J136").Value = 7066.3335
$ws.Range("L136").Value = 21199.0005
$ws.Range("N136").Value = -26299.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H134").Value = 1822.1428
$ws.Range("I134").Value = 1292.5
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 3877.5
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -1342.5
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 532.1111
$ws.Range("I5").Value = 298.16666
$ws.Range("K5").Value = 298.16666
$ws.Range("M5").Value = -186.16666
$ws.Range("H25").Value = 320
$ws.Range("I25").Value = 320
$ws.Range("K25").Value = 320
$ws.Range("M25").Value = -146
$ws.Range("H31").Value = 1724.75
$ws.Range("I31").Value = 1499.5
$ws.Range("J31").Value = 1950
$ws.Range("K31").Value = 1499.5
$ws.Range("L31").Value = 1950
$ws.Range("M31").Value = -1204.5
$ws.Range("N31").Value = -2540
$ws.Range("H34").Value = 1724.75
$ws.Range("I34").Value = 1499.5
$ws.Range("J34").Value = 1950
$ws.Range("K34").Value = 1499.5
$ws.Range("L34").Value = 1950
$ws.Range("M34").Value = -1297.5
$ws.Range("N34").Value = -2354
$ws.Range("H35").Value = 1444
$ws.Range("I35").Value = 1392.5
$ws.Range("J35").Value = 1650
$ws.Range("K35").Value = 1392.5
$ws.Range("L35").Value = 1650
$ws.Range("M35").Value = -1098.5
$ws.Range("N35").Value = -2238
$ws.Range("H58").Value = 2553.4
$ws.Range("I58").Value = 2526.111
$ws.Range("J58").Value = 2799
$ws.Range("K58").Value = 2526.111
$ws.Range("L58").Value = 2799
$ws.Range("M58").Value = -2323.111
$ws.Range("N58").Value = -3205
$ws.Range("H132").Value = 3921.5
$ws.Range("I132").Value = 4295.3335
$ws.Range("K132").Value = 12886.0005
$ws.Range("M132").Value = -10356.0005
$ws.Range("H136").Value = 2553.4
$ws.Range("I136").Value = 2526.111
$ws.Range("J136").Value = 2799
$ws.Range("K136").Value = 7578.333
$ws.Range("L136").Value = 8397
$ws.Range("M136").Value = -5028.333
$ws.Range("N136").Value = -13497

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 6166.6665
$ws.Range("I9").Value = 6166.6665
$ws.Range("K9").Value = 18499.9995
$ws.Range("M9").Value = -18275.9995
$ws.Range("H37").Value = 59982
$ws.Range("J37").Value = 59982
$ws.Range("L37").Value = 179946
$ws.Range("N37").Value = -180170
$ws.Range("H86").Value = 2669
$ws.Range("I86").Value = 410.16666
$ws.Range("J86").Value = 7186.6665
$ws.Range("K86").Value = 1230.49998
$ws.Range("L86").Value = 21559.9995
$ws.Range("M86").Value = -44.49998000000005
$ws.Range("N86").Value = -23931.9995
$ws.Range("H89").Value = 2669
$ws.Range("I89").Value = 410.16666
$ws.Range("J89").Value = 7186.6665
$ws.Range("K89").Value = 3691.49994
$ws.Range("L89").Value = 64679.9985
$ws.Range("M89").Value = 2236.50006
$ws.Range("N89").Value = -76535.9985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 327.2
$ws.Range("I97").Value = 286.8889
$ws.Range("K97").Value = 286.8889
$ws.Range("M97").Value = 209.1111
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("N113").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5912.778
$ws.Range("I7").Value = 3137.9167
$ws.Range("K7").Value = 3137.9167
$ws.Range("M7").Value = -3025.9167
$ws.Range("H22").Value = 3747
$ws.Range("I22").Value = 3796.6
$ws.Range("J22").Value = 3499
$ws.Range("K22").Value = 3796.6
$ws.Range("L22").Value = 3499
$ws.Range("M22").Value = -3501.6
$ws.Range("N22").Value = -4089
$ws.Range("H27").Value = 3747
$ws.Range("I27").Value = 3796.6
$ws.Range("J27").Value = 3499
$ws.Range("K27").Value = 3796.6
$ws.Range("L27").Value = 3499
$ws.Range("M27").Value = -3689.6
$ws.Range("N27").Value = -3713
$ws.Range("H46").Value = 38357.145
$ws.Range("I46").Value = 73214
$ws.Range("J46").Value = 3500.2856
$ws.Range("K46").Value = 73214
$ws.Range("L46").Value = 3500.2856
$ws.Range("M46").Value = -73026
$ws.Range("N46").Value = -3876.2856
$ws.Range("H61").Value = 1342.125
$ws.Range("I61").Value = 720.2857
$ws.Range("K61").Value = 720.2857
$ws.Range("M61").Value = -518.2857
$ws.Range("H93").Value = 17916.334
$ws.Range("I93").Value = 20698.8
$ws.Range("J93").Value = 4004
$ws.Range("K93").Value = 20698.8
$ws.Range("L93").Value = 4004
$ws.Range("M93").Value = -19450.8
$ws.Range("N93").Value = -6500
$ws.Range("H94").Value = 40080
$ws.Range("J94").Value = 40080
$ws.Range("L94").Value = 40080
$ws.Range("N94").Value = -41432
$ws.Range("H113").Value = 1342.125
$ws.Range("I113").Value = 720.2857
$ws.Range("K113").Value = 720.2857
$ws.Range("M113").Value = 1449.7143
$ws.Range("H122").Value = 5869.1665
$ws.Range("I122").Value = 5447.393
$ws.Range("J122").Value = 6712.7144
$ws.Range("K122").Value = 16342.179
$ws.Range("L122").Value = 20138.1432
$ws.Range("M122").Value = -13892.179
$ws.Range("N122").Value = -25038.1432
$ws.Range("H126").Value = 5912.778
$ws.Range("I126").Value = 3137.9167
$ws.Range("K126").Value = 9413.750100000001
$ws.Range("M126").Value = -6943.750100000001
$ws.Range("H132").Value = 2529.1428
$ws.Range("I132").Value = 2099.75
$ws.Range("K132").Value = 6299.25
$ws.Range("M132").Value = -3769.25
$ws.Range("H136").Value = 5999.5
$ws.Range("I136").Value = 5999
$ws.Range("K136").Value = 17997
$ws.Range("M136").Value = -15447

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7516.8335
$ws.Range("I132").Value = 8340.200000000001
$ws.Range("J132").Value = 3400
$ws.Range("K132").Value = 25020.6
$ws.Range("L132").Value = 10200
$ws.Range("M132").Value = -22490.6
$ws.Range("N132").Value = -15260
$ws.Range("H136").Value = 2804.0667
$ws.Range("J136").Value = 2418.1428
$ws.Range("L136").Value = 7254.428400000001
$ws.Range("N136").Value = -12354.4284
